# Memperbarui sistem CRUD Pengguna
# Update the attendance (presensi) recap: move/adjust entries for
# Jumat 12 Mei, Jumat 19 Mei, Minggu 21 Mei and Rabu 24 Mei 2023,
# and refresh the "Hadir" totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 - Jumat, 12 Mei 2023: clear the attendance entry (no longer recorded)
$ws.Range("B13").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("D13").ClearContents()
$ws.Range("E13").ClearContents()

# Row 20 - Jumat, 19 Mei 2023: update the check-in time, clear check-out time,
# keep status, update the remark
$ws.Range("B20").Value = "09:26:41"
$ws.Range("C20").ClearContents()
$ws.Range("D20").Value = "Hadir"
$ws.Range("E20").Value = "94,975 kilometer, TERLAMBAT 2 jam 12 menit"

# Row 22 - Minggu, 21 Mei 2023: was a "Libur Akhir Pekan" weekend row,
# now has a full attendance entry
$ws.Range("B22").Value = "21:42:21"
$ws.Range("C22").Value = "21:42:54"
$ws.Range("D22").Value = "Hadir"
$ws.Range("E22").Value = "34,744 kilometer, TERLAMBAT 14 jam 28 menit"

# Row 25 - Rabu, 24 Mei 2023: add a new attendance entry
$ws.Range("B25").Value = "17:42:33"
$ws.Range("D25").Value = "Hadir"
$ws.Range("E25").Value = "14,626 kilometer, TERLAMBAT 10 jam 28 menit"

# Update the "Hadir" (present) and overall totals from 2 to 3
$ws.Range("B34").Value = 3
$ws.Range("B37").Value = 3
